$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Banker's cheque" / customer creation rows appended to the bottom of
# the UnAuth_Customers list. The source values are text (Customer_ID /
# similar numeric-looking codes stored as strings, matching how they were
# typed in as text in the original sheet), so we force Text formatting
# before entering them, then drop the style back to Normal so the cells
# keep the default style (only their string typing persists).

$newRows = @(
    @{ Row = 100; A = "118448"; B = "17704491"; C = "1005" },
    @{ Row = 101; A = "118448"; B = "17704492"; C = "1005" }
)

foreach ($r in $newRows) {
    $rowRange = $ws.Range("A$($r.Row):C$($r.Row)")
    $rowRange.NumberFormat = "@"

    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C

    $rowRange.Style = "Normal"
}
